$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Get the existing table (ListObject) and add a new row at the end so that
# Excel automatically extends the table range and dimension.
$tbl = $ws.ListObjects.Item("Table1")
$newRow = $tbl.ListRows.Add()
$r = $newRow.Range.Row
$prevRow = $r - 1

# Copy formatting (styles + row height) from the row above, matching how the
# new row visually continues the existing table formatting.
$srcRange = $ws.Range("D" + $prevRow + ":J" + $prevRow)
$dstRange = $ws.Range("D" + $r + ":J" + $r)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item($r).RowHeight = $ws.Rows.Item($prevRow).RowHeight

# Fill in the new data row.
$ws.Cells.Item($r, 4).Value = "25/7/2029"
$ws.Cells.Item($r, 5).Value = 380
$ws.Cells.Item($r, 6).Value = 950
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 1012
$ws.Cells.Item($r, 10).Value = "N/A"

# Restore the selection state captured after the edit.
$ws.Range("E81").Select()
